$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5515
$ws1.Range("F4").Value = 2
$ws1.Range("F5").Value = 352
$ws1.Range("F6").Value = 46
$ws1.Range("C9").Value = "赣州·十万伏特-第七届青年文化综合展览会"
$ws1.Range("F9").Value = 21
$ws1.Range("F10").Value = 55
$ws1.Range("F11").Value = 128
$ws1.Range("F12").Value = 137
$ws1.Range("F13").Value = 335
$ws1.Range("F14").Value = 429
$ws1.Range("F15").Value = 3021
$ws1.Range("F16").Value = 1
$ws1.Range("F17").Value = 167
$ws1.Range("F18").Value = 1636

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5515
$ws4.Range("F4").Value = 2
$ws4.Range("F5").Value = 352
$ws4.Range("F6").Value = 46
$ws4.Range("C10").Value = "赣州·十万伏特-第七届青年文化综合展览会"
$ws4.Range("F10").Value = 21
$ws4.Range("F11").Value = 55
$ws4.Range("F12").Value = 128
$ws4.Range("F13").Value = 137
$ws4.Range("F14").Value = 335
$ws4.Range("F15").Value = 429
$ws4.Range("F16").Value = 3021
$ws4.Range("F17").Value = 1
$ws4.Range("F18").Value = 167
$ws4.Range("F19").Value = 1636
